# Generate Report for Handback
#
# The localization-status report is regenerated: the "Status" column moves
# from "Ready for handoff" to "Handed back: in sync with en-US" and each
# language sheet's "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns get populated with the handback
# results (a hyperlinked target file name, the generated xliff file name,
# and the handback timestamp for that language).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns (E, F) for both rows ---
$overview.Range("E2").Value2 = $newStatus
$overview.Range("F2").Value2 = $newStatus
$overview.Range("E3").Value2 = $newStatus
$overview.Range("F3").Value2 = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.166666666666664
$overview.Columns.Item(6).ColumnWidth = 29.166666666666664

# --- zh-cn sheet ---
$zhcn.Range("C2").Value2 = $newStatus
$zhcn.Range("C3").Value2 = $newStatus

$zhcnXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcnHandbackTime = "2016-09-06 00:42:20"

$zhcn.Range("J2").Value2 = $zhcnXlf
$zhcn.Range("K2").Value2 = $zhcnHandbackTime
$zhcn.Range("J3").Value2 = $zhcnXlf
$zhcn.Range("K3").Value2 = $zhcnHandbackTime

$aTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c0057c7329153fc5e7f4b7cca513e8f9e56589f/e2e/a.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $aTarget, "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $aTarget, "", "", "a.md")

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# --- de-de sheet ---
$dede.Range("C2").Value2 = $newStatus
$dede.Range("C3").Value2 = $newStatus

$dedeXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dedeHandbackTime = "2016-09-06 00:42:27"

$dede.Range("J2").Value2 = $dedeXlf
$dede.Range("K2").Value2 = $dedeHandbackTime
$dede.Range("J3").Value2 = $dedeXlf
$dede.Range("K3").Value2 = $dedeHandbackTime

$dede.Hyperlinks.Add($dede.Range("I2"), $aTarget, "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("I3"), $aTarget, "", "", "a.md")

$dede.Columns.Item(3).ColumnWidth = 29.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
